$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 47.875
$ws.Range("I6").Value = 40.42857
$ws.Range("K6").Value = 121.28571
$ws.Range("M6").Value = -9.285709999999995
$ws.Range("H8").Value = 368.57144
$ws.Range("I8").Value = 368.57144
$ws.Range("K8").Value = 1105.71432
$ws.Range("M8").Value = -966.71432
$ws.Range("H15").Value = 1347.0944
$ws.Range("I15").Value = 1347.0944
$ws.Range("K15").Value = 4041.2832
$ws.Range("M15").Value = -3872.2832
$ws.Range("H62").Value = 380891.12
$ws.Range("I62").Value = 9814
$ws.Range("J62").Value = 617031.0600000001
$ws.Range("K62").Value = 9814
$ws.Range("L62").Value = 617031.0600000001
$ws.Range("M62").Value = -9190
$ws.Range("N62").Value = -618279.0600000001
$ws.Range("H65").Value = 380891.12
$ws.Range("I65").Value = 9814
$ws.Range("J65").Value = 617031.0600000001
$ws.Range("K65").Value = 49070
$ws.Range("L65").Value = 3085155.3
$ws.Range("M65").Value = -45950
$ws.Range("N65").Value = -3091395.3
$ws.Range("H113").Value = 11508.23
$ws.Range("I113").Value = 19985.5
$ws.Range("J113").Value = 4242
$ws.Range("K113").Value = 19985.5
$ws.Range("L113").Value = 4242
$ws.Range("M113").Value = -16731.5
$ws.Range("N113").Value = -10750
$ws.Range("H116").Value = 1522558.8
$ws.Range("I116").Value = 7974.0835
$ws.Range("J116").Value = 3340060.5
$ws.Range("K116").Value = 7974.0835
$ws.Range("L116").Value = 3340060.5
$ws.Range("M116").Value = -4532.0835
$ws.Range("N116").Value = -3346944.5
$ws.Range("H132").Value = 1868.5
$ws.Range("I132").Value = 1363.2069
$ws.Range("K132").Value = 4089.620699999999
$ws.Range("M132").Value = -1559.620699999999
$ws.Range("H135").Value = 1310.3
$ws.Range("I135").Value = 1380.125
$ws.Range("K135").Value = 12421.125
$ws.Range("M135").Value = -9886.125

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1087
$ws.Range("I2").Value = 811.5625
$ws.Range("K2").Value = 811.5625
$ws.Range("M2").Value = -698.5625
$ws.Range("H19").Value = 651
$ws.Range("I19").Value = 651
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 651
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -422
$ws.Range("N19").ClearContents()
$ws.Range("H26").Value = 8836
$ws.Range("I26").Value = 6250
$ws.Range("J26").Value = 14008
$ws.Range("K26").Value = 6250
$ws.Range("L26").Value = 14008
$ws.Range("M26").Value = -5920
$ws.Range("N26").Value = -14668
$ws.Range("H32").Value = 9059.263000000001
$ws.Range("I32").Value = 4989.098
$ws.Range("K32").Value = 4989.098
$ws.Range("M32").Value = -4702.098
$ws.Range("H38").Value = 50647.5
$ws.Range("I38").Value = 100000
$ws.Range("J38").Value = 1295
$ws.Range("K38").Value = 100000
$ws.Range("L38").Value = 1295
$ws.Range("M38").Value = -99533
$ws.Range("N38").Value = -2229
$ws.Range("H45").Value = 2456.3076
$ws.Range("I45").Value = 1893.3
$ws.Range("K45").Value = 1893.3
$ws.Range("M45").Value = -1516.3
$ws.Range("H61").Value = 2411.2666
$ws.Range("I61").Value = 1980.8334
$ws.Range("K61").Value = 1980.8334
$ws.Range("M61").Value = -1768.8334
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H116").Value = 1087
$ws.Range("I116").Value = 811.5625
$ws.Range("K116").Value = 811.5625
$ws.Range("M116").Value = 1482.4375
$ws.Range("H136").Value = 2411.2666
$ws.Range("I136").Value = 1980.8334
$ws.Range("K136").Value = 5942.5002
$ws.Range("M136").Value = -3392.5002

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1087
$ws.Range("I3").Value = 811.5625
$ws.Range("K3").Value = 811.5625
$ws.Range("M3").Value = -697.5625
$ws.Range("H12").Value = 808.7143
$ws.Range("I12").Value = 776.8333
$ws.Range("J12").Value = 1000
$ws.Range("K12").Value = 776.8333
$ws.Range("L12").Value = 1000
$ws.Range("M12").Value = -608.8333
$ws.Range("N12").Value = -1336
$ws.Range("H40").Value = 14800
$ws.Range("J40").Value = 14800
$ws.Range("L40").Value = 14800
$ws.Range("N40").Value = -15330
$ws.Range("H86").Value = 2129.2666
$ws.Range("I86").Value = 2244.875
$ws.Range("J86").Value = 1997.1428
$ws.Range("K86").Value = 2244.875
$ws.Range("L86").Value = 1997.1428
$ws.Range("M86").Value = -1121.875
$ws.Range("N86").Value = -4243.1428
$ws.Range("H89").Value = 2129.2666
$ws.Range("I89").Value = 2244.875
$ws.Range("J89").Value = 1997.1428
$ws.Range("K89").Value = 11224.375
$ws.Range("L89").Value = 9985.714
$ws.Range("M89").Value = -5608.375
$ws.Range("N89").Value = -21217.714
$ws.Range("H94").Value = 456.3
$ws.Range("I94").Value = 428.4
$ws.Range("J94").Value = 540
$ws.Range("K94").Value = 428.4
$ws.Range("L94").Value = 540
$ws.Range("M94").Value = 22.60000000000002
$ws.Range("N94").Value = -1442
$ws.Range("H99").Value = 8087505
$ws.Range("I99").Value = 525005
$ws.Range("K99").Value = 525005
$ws.Range("M99").Value = -523507
$ws.Range("H107").Value = 4071.2942
$ws.Range("I107").Value = 3100.6365
$ws.Range("J107").Value = 5850.8335
$ws.Range("K107").Value = 3100.6365
$ws.Range("L107").Value = 5850.8335
$ws.Range("M107").Value = -1180.6365
$ws.Range("N107").Value = -9690.833500000001
$ws.Range("H134").Value = 1416.7843
$ws.Range("I134").Value = 1090.1333
$ws.Range("K134").Value = 3270.3999
$ws.Range("M134").Value = -735.3998999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 20599.8
$ws.Range("I15").Value = 1499.5
$ws.Range("J15").Value = 33333.332
$ws.Range("K15").Value = 1499.5
$ws.Range("L15").Value = 33333.332
$ws.Range("M15").Value = -1329.5
$ws.Range("N15").Value = -33673.332
$ws.Range("H16").Value = 1824.2858
$ws.Range("I16").Value = 1563.8823
$ws.Range("K16").Value = 1563.8823
$ws.Range("M16").Value = -1276.8823
$ws.Range("H58").Value = 2357
$ws.Range("I58").Value = 2265.1428
$ws.Range("J58").Value = 3000
$ws.Range("K58").Value = 2265.1428
$ws.Range("L58").Value = 3000
$ws.Range("M58").Value = -2062.1428
$ws.Range("N58").Value = -3406
$ws.Range("H113").Value = 1824.2858
$ws.Range("I113").Value = 1563.8823
$ws.Range("K113").Value = 1563.8823
$ws.Range("M113").Value = 606.1177
$ws.Range("H136").Value = 2357
$ws.Range("I136").Value = 2265.1428
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 6795.428400000001
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -4245.428400000001
$ws.Range("N136").Value = -14100

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 5291.5
$ws.Range("I7").Value = 70
$ws.Range("K7").Value = 210
$ws.Range("M7").Value = -98

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 27000
$ws.Range("J32").Value = 27000
$ws.Range("L32").Value = 27000
$ws.Range("N32").Value = -27592
$ws.Range("H51").Value = 50000
$ws.Range("J51").Value = 50000
$ws.Range("L51").Value = 50000
$ws.Range("N51").Value = -51018
$ws.Range("H122").Value = 2773.2856
$ws.Range("I122").Value = 2034.5454
$ws.Range("J122").Value = 5482
$ws.Range("K122").Value = 6103.6362
$ws.Range("L122").Value = 16446
$ws.Range("M122").Value = -3653.6362
$ws.Range("N122").Value = -21346

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 22521.217
$ws.Range("I7").Value = 11761.333
$ws.Range("K7").Value = 11761.333
$ws.Range("M7").Value = -11649.333
$ws.Range("H16").Value = 3666.2222
$ws.Range("I16").Value = 4285.2856
$ws.Range("K16").Value = 4285.2856
$ws.Range("M16").Value = -4115.2856
$ws.Range("H36").Value = 87499.5
$ws.Range("J36").Value = 87499.5
$ws.Range("L36").Value = 87499.5
$ws.Range("N36").Value = -88623.5
$ws.Range("H40").Value = 7939079.5
$ws.Range("I40").Value = 3000.1667
$ws.Range("K40").Value = 3000.1667
$ws.Range("M40").Value = -2864.1667
$ws.Range("H61").Value = 568.4
$ws.Range("I61").Value = 568.4
$ws.Range("K61").Value = 568.4
$ws.Range("M61").Value = -366.4
$ws.Range("H113").Value = 568.4
$ws.Range("I113").Value = 568.4
$ws.Range("K113").Value = 568.4
$ws.Range("M113").Value = 1601.6
$ws.Range("H122").Value = 10024004
$ws.Range("I122").Value = 29067.188
$ws.Range("J122").Value = 50003750
$ws.Range("K122").Value = 87201.564
$ws.Range("L122").Value = 150011250
$ws.Range("M122").Value = -84751.564
$ws.Range("N122").Value = -150016150
$ws.Range("H126").Value = 22521.217
$ws.Range("I126").Value = 11761.333
$ws.Range("K126").Value = 35283.999
$ws.Range("M126").Value = -32813.999
$ws.Range("H136").Value = 5429.773
$ws.Range("I136").Value = 6897.1
$ws.Range("K136").Value = 20691.3
$ws.Range("M136").Value = -18141.3

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1950.2858
$ws.Range("I107").Value = 1453.7646
$ws.Range("K107").Value = 4361.293799999999
$ws.Range("M107").Value = -2441.293799999999
$ws.Range("H122").Value = 3151
$ws.Range("I122").Value = 2441.3333
$ws.Range("K122").Value = 7323.999899999999
$ws.Range("M122").Value = -4873.999899999999
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H141").Value = 79071.336
$ws.Range("J141").Value = 79071.336
$ws.Range("L141").Value = 79071.336
$ws.Range("N141").Value = -89431.336
